$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update election result totals in row 2 (GUARDA / PINHEL), columns H:AA
$ws.Range("H2").Value  = 92
$ws.Range("I2").Value  = 244
$ws.Range("J2").Value  = 1024
$ws.Range("K2").Value  = 5
$ws.Range("L2").Value  = 280
$ws.Range("M2").Value  = 18
$ws.Range("N2").Value  = 195
$ws.Range("O2").Value  = 1
$ws.Range("P2").Value  = 1
$ws.Range("Q2").Value  = 0
$ws.Range("R2").Value  = 12
$ws.Range("S2").Value  = 130
$ws.Range("T2").Value  = 164
$ws.Range("U2").Value  = 15
$ws.Range("V2").Value  = 1683
$ws.Range("W2").Value  = 1
$ws.Range("X2").Value  = 1651
$ws.Range("Y2").Value  = 0
$ws.Range("Z2").Value  = 31
$ws.Range("AA2").Value = 10
